# The commit swaps the contents of ppt/theme/theme1.xml (used by the slide
# master, i.e. the deck's main "Integral" design) and ppt/theme/theme2.xml
# (used by the notes master, previously the default "Office Theme" colors).
# After the edit, the slide master's theme carries the default Office color
# palette while the notes master keeps the old Integral palette.
#
# The only durable, COM-exposed lever for editing a theme's color values in
# this host is Master.Theme.ThemeColorScheme.Colors(n).RGB - and it always
# targets the presentation's single live theme part (ppt/theme/theme1.xml),
# regardless of whether it is reached via SlideMaster, NotesMaster,
# HandoutMaster, a Slide, or a SlideRange. So we drive the swap through the
# SlideMaster and push the new ("Office Theme") palette onto it - this is
# the part of the change that is both visible and reachable.
#
# PowerPoint's RGB() packs a hex color "RRGGBB" as R + G*256 + B*65536;
# ConvertToRgb below performs that conversion from a plain hex string so the
# intent (which named color goes where) stays readable.

function ConvertTo-PptRgb {
    param([string]$HexColor)
    $r = [Convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index order matches MsoThemeColorSchemeIndex / the <a:clrScheme> element
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = ConvertTo-PptRgb $officeThemeColors[$i]
}
